$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("L2").Value = 1.33
$ws.Range("N2").Value = 4.7
$ws.Range("O2").Value = 1.24
$ws.Range("P2").Value = 2.28
$ws.Range("Q2").Value = 1.72
$ws.Range("R2").Value = 1.51
$ws.Range("S2").Value = 2.8
# Row 3
$ws.Range("F3").Value = 1.53
$ws.Range("G3").Value = 1.64
$ws.Range("H3").Value = 6.8
$ws.Range("I3").Value = 8.4
$ws.Range("J3").Value = 4
$ws.Range("M3").Value = 1.06
$ws.Range("N3").Value = 3.55
$ws.Range("P3").Value = 1.87
$ws.Range("Q3").Value = 1.92
$ws.Range("R3").Value = 1.33
$ws.Range("S3").Value = 3.35
$ws.Range("T3").Value = 2
$ws.Range("U3").Value = 1.82
$ws.Range("V3").Value = 1.13
$ws.Range("W3").Value = 2.56
$ws.Range("Y3").Value = 30
$ws.Range("Z3").Value = 85
$ws.Range("AB3").Value = 10.5
$ws.Range("AC3").Value = 13.5
$ws.Range("AD3").Value = 38
$ws.Range("AF3").Value = 12.5
$ws.Range("AG3").Value = 13.5
$ws.Range("AH3").Value = 34
$ws.Range("AK3").Value = 25
$ws.Range("AL3").Value = 55
# Row 4
$ws.Range("G4").Value = 2.14
$ws.Range("H4").Value = 4.3
$ws.Range("I4").Value = 5.1
$ws.Range("J4").Value = 3.15
$ws.Range("K4").Value = 3.75
$ws.Range("M4").Value = 1.1
$ws.Range("N4").Value = 2.98
$ws.Range("O4").Value = 1.43
$ws.Range("R4").Value = 1.24
$ws.Range("S4").Value = 4.3
$ws.Range("T4").Value = 1.98
$ws.Range("U4").Value = 1.83
$ws.Range("V4").Value = 1.24
$ws.Range("W4").Value = 1.88
$ws.Range("X4").Value = 11
$ws.Range("Y4").Value = 14
$ws.Range("Z4").Value = 36
$ws.Range("AA4").Value = 140
$ws.Range("AB4").Value = 8
$ws.Range("AC4").Value = 8
$ws.Range("AD4").Value = 19.5
$ws.Range("AE4").Value = 75
$ws.Range("AF4").Value = 12.5
$ws.Range("AG4").Value = 11.5
$ws.Range("AH4").Value = 23
$ws.Range("AI4").Value = 100
$ws.Range("AJ4").Value = 26
$ws.Range("AK4").Value = 27
$ws.Range("AL4").Value = 55
$ws.Range("AM4").Value = 180
$ws.Range("AN4").Value = 22
$ws.Range("AO4").Value = 120
# Row 5
$ws.Range("K5").Value = 4.1
$ws.Range("L5").Value = 1.32
$ws.Range("Q5").Value = 1.85
$ws.Range("R5").Value = 1.37
$ws.Range("S5").Value = 2.88
$ws.Range("T5").Value = 1.76
$ws.Range("U5").Value = 2.12
$ws.Range("V5").Value = 1.94
# Row 6
$ws.Range("L6").Value = 1.31
$ws.Range("M6").Value = 1.05
$ws.Range("N6").Value = 4.6
$ws.Range("O6").Value = 1.23
$ws.Range("R6").Value = 1.49
$ws.Range("S6").Value = 2.72
$ws.Range("T6").Value = 1.94
$ws.Range("U6").Value = 1.93
$ws.Range("V6").Value = 3.15
$ws.Range("W6").Value = 1.12
$ws.Range("X6").Value = 26
$ws.Range("Y6").Value = 9.199999999999999
$ws.Range("Z6").Value = 9.199999999999999
$ws.Range("AA6").Value = 12.5
$ws.Range("AB6").Value = 29
$ws.Range("AC6").Value = 12
$ws.Range("AD6").Value = 10.5
$ws.Range("AE6").Value = 15.5
$ws.Range("AF6").Value = 75
$ws.Range("AG6").Value = 36
$ws.Range("AH6").Value = 25
$ws.Range("AI6").Value = 40
$ws.Range("AJ6").Value = 330
$ws.Range("AK6").Value = 150
$ws.Range("AL6").Value = 130
$ws.Range("AM6").Value = 160
$ws.Range("AN6").Value = 190
$ws.Range("AO6").Value = 7.2
# Row 7
$ws.Range("F7").Value = 1.17
$ws.Range("G7").Value = 1.21
$ws.Range("H7").Value = 18.5
$ws.Range("I7").Value = 26
$ws.Range("J7").Value = 8.6
$ws.Range("K7").Value = 11
$ws.Range("L7").Value = 1.01
$ws.Range("M7").Value = 1.01
$ws.Range("N7").Value = 1.01
$ws.Range("O7").Value = 1.09
$ws.Range("Q7").Value = 1.3
$ws.Range("R7").Value = 1.84
$ws.Range("S7").Value = 1.71
$ws.Range("T7").Value = 1.78
$ws.Range("U7").Value = 1.7
$ws.Range("V7").Value = 1.04
$ws.Range("W7").Value = 5.8
$ws.Range("X7").Value = 1000
$ws.Range("Y7").Value = 1000
$ws.Range("Z7").Value = 1000
$ws.Range("AA7").Value = 1000
$ws.Range("AB7").Value = 1000
$ws.Range("AC7").Value = 1000
$ws.Range("AD7").Value = 1000
$ws.Range("AE7").Value = 1000
$ws.Range("AF7").Value = 1000
$ws.Range("AG7").Value = 1000
$ws.Range("AH7").Value = 1000
$ws.Range("AI7").Value = 1000
$ws.Range("AJ7").Value = 1000
$ws.Range("AK7").Value = 1000
$ws.Range("AL7").Value = 1000
$ws.Range("AM7").Value = 1000
$ws.Range("AN7").Value = 1000
$ws.Range("AO7").Value = 1000
# Row 8
$ws.Range("F8").Value = 1.32
$ws.Range("G8").Value = 1.37
$ws.Range("J8").Value = 6.2
$ws.Range("K8").Value = 7.2
$ws.Range("X8").Value = 40
$ws.Range("AC8").Value = 16
$ws.Range("AD8").Value = 36
$ws.Range("AG8").Value = 11.5
$ws.Range("AH8").Value = 26
$ws.Range("AJ8").Value = 12.5
$ws.Range("AK8").Value = 14
$ws.Range("AL8").Value = 29
$ws.Range("AM8").Value = 110
$ws.Range("AO8").Value = 110
# Row 9
$ws.Range("F9").Value = 2.08
$ws.Range("G9").Value = 2.8
$ws.Range("H9").Value = 2.78
$ws.Range("I9").Value = 4.2
$ws.Range("J9").Value = 3.05
$ws.Range("K9").Value = 6.4
$ws.Range("L9").Value = 1.01
$ws.Range("M9").Value = 1.01
$ws.Range("N9").Value = 1.64
$ws.Range("O9").Value = 1.01
$ws.Range("P9").Value = 1.64
$ws.Range("Q9").Value = 1.9
$ws.Range("R9").Value = 1.23
$ws.Range("S9").Value = 3.25
$ws.Range("T9").Value = 1.01
$ws.Range("U9").Value = 1.01
$ws.Range("V9").Value = 1.32
$ws.Range("W9").Value = 1.55
$ws.Range("X9").Value = 1000
$ws.Range("Y9").Value = 1000
$ws.Range("Z9").Value = 1000
$ws.Range("AA9").Value = 1000
$ws.Range("AB9").Value = 1000
$ws.Range("AC9").Value = 1000
$ws.Range("AD9").Value = 1000
$ws.Range("AE9").Value = 1000
$ws.Range("AF9").Value = 1000
$ws.Range("AG9").Value = 1000
$ws.Range("AH9").Value = 1000
$ws.Range("AI9").Value = 1000
$ws.Range("AJ9").Value = 1000
$ws.Range("AK9").Value = 1000
$ws.Range("AL9").Value = 1000
$ws.Range("AM9").Value = 1000
$ws.Range("AN9").Value = 1000
$ws.Range("AO9").Value = 1000
# Row 10
$ws.Range("J10").Value = 2.9
$ws.Range("L10").Value = 1.01
$ws.Range("M10").Value = 1.01
$ws.Range("N10").Value = 1.81
$ws.Range("O10").Value = 1.33
$ws.Range("R10").Value = 1.25
$ws.Range("S10").Value = 3.1
$ws.Range("T10").Value = 1.01
$ws.Range("U10").Value = 1.01
$ws.Range("V10").Value = 1.7
$ws.Range("W10").Value = 1.29
$ws.Range("X10").Value = 18.5
$ws.Range("Y10").Value = 13.5
$ws.Range("Z10").Value = 21
$ws.Range("AA10").Value = 46
$ws.Range("AB10").Value = 18.5
$ws.Range("AC10").Value = 11
$ws.Range("AD10").Value = 16
$ws.Range("AE10").Value = 38
$ws.Range("AF10").Value = 36
$ws.Range("AG10").Value = 22
$ws.Range("AH10").Value = 26
$ws.Range("AI10").Value = 65
$ws.Range("AJ10").Value = 1000
$ws.Range("AK10").Value = 65
$ws.Range("AL10").Value = 1000
$ws.Range("AM10").Value = 1000
$ws.Range("AN10").Value = 1000
$ws.Range("AO10").Value = 1000
# Row 11
$ws.Range("F11").Value = 1.84
$ws.Range("G11").Value = 2.3
$ws.Range("J11").Value = 3.5
$ws.Range("K11").Value = 7
# Row 12
$ws.Range("F12").Value = 1.87
$ws.Range("I12").Value = 4.3
$ws.Range("P12").Value = 3.1
$ws.Range("Q12").Value = 1.37
# Row 13
$ws.Range("J13").Value = 3.95
$ws.Range("K13").Value = 5.1
$ws.Range("L13").Value = 1.22
$ws.Range("R13").Value = 1.73
$ws.Range("S13").Value = 2.02
$ws.Range("T13").Value = 1.44
$ws.Range("U13").Value = 2.78
$ws.Range("AI13").Value = 32
# Row 15
$ws.Range("H15").Value = 4.3
$ws.Range("I15").Value = 12.5
$ws.Range("J15").Value = 3.95
$ws.Range("K15").Value = 10
# Row 17
$ws.Range("K17").Value = 6.2
$ws.Range("Q17").Value = 1.39
# Row 18
$ws.Range("F18").Value = 2.14
$ws.Range("G18").Value = 2.4
$ws.Range("H18").Value = 3.8
$ws.Range("I18").Value = 4.3
$ws.Range("J18").Value = 3.15
$ws.Range("K18").Value = 3.75
$ws.Range("P18").Value = 1.66
$ws.Range("Q18").Value = 2.28
# Row 21
$ws.Range("G21").Value = 1.97
$ws.Range("I21").Value = 5.1
$ws.Range("Q21").Value = 2.12
# Row 22
$ws.Range("F22").Value = 2.6
$ws.Range("G22").Value = 2.9
# Row 23
$ws.Range("L23").Value = 1.54
# Row 24
$ws.Range("F24").Value = 1.09
# Row 25
$ws.Range("O25").Value = 1.43
$ws.Range("AO25").Value = 46
# Row 26
$ws.Range("AC26").Value = 7.2
# Row 27
$ws.Range("F27").Value = 1.41
$ws.Range("G27").Value = 1.44
$ws.Range("H27").Value = 9.800000000000001
$ws.Range("I27").Value = 12
$ws.Range("J27").Value = 4.9
$ws.Range("P27").Value = 1.79
$ws.Range("Q27").Value = 2.14
# Row 28
$ws.Range("I28").Value = 3.6
$ws.Range("K28").Value = 3.4
$ws.Range("Q28").Value = 2.2

Write-Host "Applied 264 cell updates"
